$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New "zeef" section header (bold, like the "circuit satisfiability problem" header)
$ws.Range("B36").Value = "zeef"
$ws.Range("B36").Font.Bold = $true

# Column headers for the new table (set C37 first so shared-string order matches)
$ws.Range("C37").Value = "threads & processors"
$ws.Range("B37").Value = "runtime"

# Data rows
$ws.Range("B38").Value = 2
$ws.Range("C38").Value = 1
$ws.Range("C39").Value = 2
$ws.Range("C40").Value = 4
$ws.Range("C41").Value = 8

# Update sheet view to match the new selection/scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B38").Select()
$excel.ActiveWindow.zoomScaleNormal = 80
